$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") holds a date serial that was bumped by one day
# (2024-11-01 -> 2024-11-02) for every data row (rows 2 through 32).
for ($row = 2; $row -le 32; $row++) {
    $ws.Cells.Item($row, 3).Value2 = 45598
}
